$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Bad Drivers" top entry (row 3) with new sample numbers ---
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 369
$ws.Range("D3").Value = 98.8

# --- Remove the second "Bad Drivers" entry (old row 4); Totals row moves up to row 4 ---
$ws.Rows.Item(4).Delete()

# Update the Totals row (now row 4) to reflect the remaining single entry
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 369

# --- Remove the first "Good Drivers" entry (AX201 23.100.0.4), now at row 12 ---
$ws.Rows.Item(12).Delete()

# --- Remove the trailing "Good Drivers" entries that are no longer reported ---
# After the previous deletions these are rows 15 through 19 (5 rows):
#   AX201 22.80.0.9, AX201 22.50.1.1, AX201 21.110.3.2, AX201 21.70.0.6, AX201 21.60.2.1
$ws.Range("A15:A19").EntireRow.Delete()

# --- Update remaining "Good Drivers" rows with the new figures ---
# Row 12: AX211 22.150.3.1 - Driver Vintage not available this week (blank numeric placeholder)
$ws.Range("E12").Value = 0

# Row 14: AX211 22.100.1.1 - Total Samples updated
$ws.Range("B14").Value = 265400
